# LCFS edits from WRI
# Applies the content-level changes described by the commit:
#  - About sheet: replace the two-sentence note about railways being
#    included (implementing ministries / National Biofuel Policy) with a
#    single sentence saying railways are excluded, and remove the now
#    unused second line.
#  - BVTStL sheet: flip which vehicle types are subject to LCFS for the
#    "rail" and "ships" rows (rail: 1,1 -> 0,0 ; ships: 0,0 -> 1,0).
#  - Update the remembered cell selections on both sheets.

$wb  = $excel.ActiveWorkbook
$wsAbout  = $wb.Worksheets.Item("About")
$wsBVTStL = $wb.Worksheets.Item("BVTStL")

# --- About sheet: update the note text about railways -------------------
# Old text (two shared-string rows):
#   "We include railways as it is listed as one of the implementing ministries"
#   "in the National Biofuel Policy document. "
# New text (single row, replacing the first; the second is cleared out):
$wsAbout.Range("A16").Value = "We therefore exclude railways."
$wsAbout.Range("A17").ClearContents()

# --- BVTStL sheet: rail is no longer subject to LCFS, ships now is ------
$wsBVTStL.Range("B5").Value = 0   # rail / passenger
$wsBVTStL.Range("C5").Value = 0   # rail / freight
$wsBVTStL.Range("B6").Value = 1   # ships / passenger
$wsBVTStL.Range("C6").Value = 0   # ships / freight (unchanged)

# --- Remembered selections ----------------------------------------------
# Select on the non-active sheet first, then the active sheet last so the
# workbook keeps "About" as the selected/visible tab.
$wsBVTStL.Range("F9").Select()
$wsAbout.Range("B38").Select()
